# Sync attendance_reports: reorder "Recorded By" (column G) value lists so
# that comma-separated entries are reversed in order (e.g. "a, System" -> "System, a").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $val = $cell.Value2
    if ($null -ne $val -and $val -is [string] -and $val.Contains(",")) {
        $parts = $val -split ","
        for ($i = 0; $i -lt $parts.Length; $i++) {
            $parts[$i] = $parts[$i].Trim()
        }

        # Only reorder lists that explicitly include "System" as one of the
        # recorder names (moves it to the front by reversing the list, since
        # it is always appended last in the original data).
        $hasSystem = $false
        foreach ($p in $parts) {
            if ($p.Equals("System")) {
                $hasSystem = $true
            }
        }

        if ($hasSystem) {
            $reversedParts = @()
            for ($i = $parts.Length - 1; $i -ge 0; $i--) {
                $reversedParts += $parts[$i]
            }

            $cell.Value2 = [string]::Join(", ", $reversedParts)
        }
    }
}
